# Adapt tests to control version
# Adds a "version" column to the settings sheet (header + value 1),
# and makes "settings" the active sheet/tab of the workbook.

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("settings")

# Add the new "version" column header and value next to form_title/form_id.
$settings.Range("C1").Value = "version"
$settings.Range("C2").Value = 1

# Make "settings" the active sheet (activeTab=2 / tabSelected on settings sheet),
# with the new cell as the active selection.
$settings.Activate()
$settings.Range("C3").Select()
